# Re-sort the stock list (Sheet1) alphabetically by QuoteSymbol (column A),
# keeping the header row (row 1) and the first data row (^GSPC, row 2) fixed
# in place, while turning on AutoFilter for the table - reproducing the
# "Data > Sort" + AutoFilter trace left behind in the workbook XML.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Turn on AutoFilter for the table range.
$ws.Range("A1:B8").AutoFilter()

# Sort rows A2:B8 (data below the fixed ^GSPC row) alphabetically by column A,
# leaving the header (row 1) out of the sort.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A8"))
$ws.Sort.SetRange($ws.Range("A1:B8"))
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Match the resulting selection left behind in the file.
$ws.Range("A3").Select()
